# Add three new test-case blocks (vacations / "Affichage" tests) to the sheet,
# mirroring the layout of the existing blocks (section title row, column
# header row, then one data row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Block 1: "Affichage d'un pilote" (rows 31-33) ----------------------
$ws.Range("A20").Copy()
$ws.Range("A31").PasteSpecial(-4122)
$ws.Range("A31").Value = "Affichage d'un pilote"

$ws.Range("A21:D21").Copy()
$ws.Range("A32:D32").PasteSpecial(-4122)
$ws.Range("A32").Value = "`"Champs`""
$ws.Range("B32").Value = "Comportement Attendu"
$ws.Range("C32").Value = "Comportement Observé"
$ws.Range("D32").Value = "Résultat"

$ws.Range("A23:D23").Copy()
$ws.Range("A33:D33").PasteSpecial(-4122)
$ws.Range("A33").Value = "DataGridViewer"
$ws.Range("B33").Value = "Pour chaque pilote répertorié dans la base de données sont affichées par ligne les informations suivantes : Id, Nom, Prénom, Aéroport d'affectation, Heures de vol"
$ws.Range("C33").Value = "Les informations sont retrouvées dans les cases correspondantes du dataGridView pour chaque pilote"
$ws.Range("D33").Value = "OK"
$ws.Range("A33:D33").RowHeight = 60

# ---- Block 2: "Affichage d'un vol" (rows 37-39) --------------------------
$ws.Range("A20").Copy()
$ws.Range("A37").PasteSpecial(-4122)
$ws.Range("A37").Value = "Affichage d'un vol"

$ws.Range("A21:D21").Copy()
$ws.Range("A38:D38").PasteSpecial(-4122)
$ws.Range("A38").Value = "`"Champs`""
$ws.Range("B38").Value = "Comportement Attendu"
$ws.Range("C38").Value = "Comportement Observé"
$ws.Range("D38").Value = "Résultat"

$ws.Range("A23:D23").Copy()
$ws.Range("A39:D39").PasteSpecial(-4122)
$ws.Range("A39").Value = "DataGridViewer"
$ws.Range("B39").Value = "Pour chaque vol répertorié dans la base de données sont affichées par ligne les informations suivantes : Nom, Ligne, Date départ, Date arrivée, Pilot n°1, Pilote n°2"
$ws.Range("C39").Value = "Les informations sont retrouvées dans les cases correspondantes du dataGridView pour chaque vol. Mais comme l'affectation d'un vol à un pilote n'est pas encore implémantée, je ne peux pas verifier si l'information est correcte pour Pilote n°1 et Pilote n°2"
$ws.Range("D39").Value = "OK"
$ws.Range("A39:D39").RowHeight = 75

# ---- Block 3: "Affichage d'une ligne" (rows 42-44) -----------------------
$ws.Range("A20").Copy()
$ws.Range("A42").PasteSpecial(-4122)
$ws.Range("A42").Value = "Affichage d'une ligne"

$ws.Range("A21:D21").Copy()
$ws.Range("A43:D43").PasteSpecial(-4122)
$ws.Range("A43").Value = "`"Champs`""
$ws.Range("B43").Value = "Comportement Attendu"
$ws.Range("C43").Value = "Comportement Observé"
$ws.Range("D43").Value = "Résultat"

$ws.Range("A23:D23").Copy()
$ws.Range("A44:D44").PasteSpecial(-4122)
$ws.Range("A44").Value = "DataGridViewer"
$ws.Range("B44").Value = "Pour chaque ligne répertoriée dans la base de données sont affichées par ligne les informations suivantes : Id, lieu de départ, lieu d'arrivée, distance"
$ws.Range("C44").Value = "Les informations sont retrouvées dans les cases correspondantes du dataGridView pour chaque ligne"
$ws.Range("D44").Value = "OK"
$ws.Range("A44:D44").RowHeight = 60

# ---- View state: scroll down to the newly added rows ---------------------
$ws.Application.ActiveWindow.ScrollRow = 32
$ws.Range("C48").Select()
